$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells we are about to rewrite to stay as Text,
# so purely-numeric-looking strings (e.g. "0.9988") are not auto-converted
# into numbers by Excel's normal value-parsing.
$ws.Range("D2:D4").NumberFormat = "@"
$ws.Range("D6:D10").NumberFormat = "@"
$ws.Range("D12:D20").NumberFormat = "@"
$ws.Range("D22:D51").NumberFormat = "@"

$ws.Range("D2").Value = '23.439.16'
$ws.Range("E2").Value = '  -1.06%  '
$ws.Range("D3").Value = '1.646.66'
$ws.Range("E3").Value = '  -0.48%  '
$ws.Range("D4").Value = '0.9988'
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("E5").Value = '  -0.31%  '
$ws.Range("D6").Value = '299.88'
$ws.Range("E6").Value = '  -1.06%  '
$ws.Range("D7").Value = '0.3798'
$ws.Range("E7").Value = '  -1.14%  '
$ws.Range("D8").Value = '50.21'
$ws.Range("E8").Value = '  -1.68%  '
$ws.Range("D9").Value = '0.3496'
$ws.Range("E9").Value = '  -3.01%  '
$ws.Range("D10").Value = '0.08069'
$ws.Range("E10").Value = '  -1.53%  '
$ws.Range("E11").Value = '  -0.94%  '
$ws.Range("D12").Value = '0.9988'
$ws.Range("E12").Value = '  -0.27%  '
$ws.Range("D13").Value = '22.04'
$ws.Range("E13").Value = '  -1.66%  '
$ws.Range("D14").Value = '6.323'
$ws.Range("E14").Value = '  -2.12%  '
$ws.Range("D15").Value = '7.276'
$ws.Range("E15").Value = '  -2.20%  '
$ws.Range("D16").Value = '0.00001214'
$ws.Range("E16").Value = '  -0.61%  '
$ws.Range("D17").Value = '1.646.38'
$ws.Range("E17").Value = '  -0.41%  '
$ws.Range("D18").Value = '94.83'
$ws.Range("E18").Value = '  -2.82%  '
$ws.Range("D19").Value = '0.06966'
$ws.Range("E19").Value = '  -1.00%  '
$ws.Range("D20").Value = '6.632'
$ws.Range("E20").Value = '  -2.03%  '
$ws.Range("E21").Value = '  -0.90%  '
$ws.Range("D22").Value = '0.9986'
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("D23").Value = '12.40'
$ws.Range("E23").Value = '  -2.18%  '
$ws.Range("D24").Value = '23.456.69'
$ws.Range("E24").Value = '  -0.96%  '
$ws.Range("D25").Value = '2.434'
$ws.Range("E25").Value = '  -1.95%  '
$ws.Range("D26").Value = '2.968'
$ws.Range("E26").Value = '  -2.03%  '
$ws.Range("D27").Value = '21.03'
$ws.Range("E27").Value = '  -1.03%  '
$ws.Range("D28").Value = '149.81'
$ws.Range("E28").Value = '  -2.75%  '
$ws.Range("D29").Value = '5.182'
$ws.Range("E29").Value = '  -1.25%  '
$ws.Range("D30").Value = '131.87'
$ws.Range("E30").Value = '  -1.63%  '
$ws.Range("D31").Value = '1.816.67'
$ws.Range("E31").Value = '  -1.10%  '
$ws.Range("D32").Value = '6.853'
$ws.Range("E32").Value = '  -3.31%  '
$ws.Range("D33").Value = '2.128'
$ws.Range("E33").Value = '  -5.58%  '
$ws.Range("D34").Value = '11.26'
$ws.Range("E34").Value = '  -6.78%  '
$ws.Range("D35").Value = '0.9896'
$ws.Range("E35").Value = '  -6.54%  '
$ws.Range("D36").Value = '0.02686'
$ws.Range("E36").Value = '  -4.31%  '
$ws.Range("D37").Value = '0.08787'
$ws.Range("E37").Value = '  -0.14%  '
$ws.Range("D38").Value = '0.2423'
$ws.Range("E38").Value = '  -3.32%  '
$ws.Range("D39").Value = '5.893'
$ws.Range("E39").Value = '  -3.08%  '
$ws.Range("D40").Value = '0.06833'
$ws.Range("E40").Value = '  -2.16%  '
$ws.Range("D41").Value = '12.78'
$ws.Range("E41").Value = '  -1.96%  '
$ws.Range("D42").Value = '0.6833'
$ws.Range("E42").Value = '  -2.22%  '
$ws.Range("D43").Value = '1.289'
$ws.Range("E43").Value = '  -3.61%  '
$ws.Range("D44").Value = '15.45'
$ws.Range("E44").Value = '  -3.36%  '
$ws.Range("D45").Value = '0.9929'
$ws.Range("E45").Value = '  -0.85%  '
$ws.Range("D46").Value = '0.6353'
$ws.Range("E46").Value = '  -2.33%  '
$ws.Range("D47").Value = '2.241'
$ws.Range("E47").Value = '  -2.64%  '
$ws.Range("D48").Value = '3.912'
$ws.Range("E48").Value = '  -1.30%  '
$ws.Range("D49").Value = '127.27'
$ws.Range("E49").Value = '  -0.71%  '
$ws.Range("D50").Value = '0.07687'
$ws.Range("E50").Value = '  -2.48%  '
$ws.Range("D51").Value = '1.228'
$ws.Range("E51").Value = '  +2.92%  '
